$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# G4/G5 share the same "Latest HO Xliff Generate Date" string (2016-08-25 16:16:16)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-25 16:17:14"
$wsOverview.Range("G5").Value = "2016-08-25 16:17:14"

# --- zh-cn sheet ---
# E4/E5 share "Priority" = ht -> mt
# H4/H5 share "Correspond Handoff Datetime" = 2016-08-25 16:16:10 -> 2016-08-25 16:17:09
# K4/K5 share "Correspond Handback DateTime" = 2016-08-25 16:16:36 -> 2016-08-25 16:17:34
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-25 16:17:09"
$wsZhCn.Range("H5").Value = "2016-08-25 16:17:09"
$wsZhCn.Range("K4").Value = "2016-08-25 16:17:34"
$wsZhCn.Range("K5").Value = "2016-08-25 16:17:34"

# --- de-de sheet ---
# E4/E5 share "Priority" = ht -> mt (same si as zh-cn E4/E5)
# H4/H5 share "Correspond Handoff Datetime" = 2016-08-25 16:16:16 -> 2016-08-25 16:17:14 (same si as Overview G4/G5)
# K4/K5 share "Correspond Handback DateTime" = 2016-08-25 16:16:44 -> 2016-08-25 16:17:42
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-25 16:17:14"
$wsDeDe.Range("H5").Value = "2016-08-25 16:17:14"
$wsDeDe.Range("K4").Value = "2016-08-25 16:17:42"
$wsDeDe.Range("K5").Value = "2016-08-25 16:17:42"
